# Product mix 5 test cases, Data tables 15 Test Cases
#
# - Deselect "Modify Transaction" tab, select "Share Other Details1" tab
#   (active tab moves from index 2 to index 3)
# - On "Share Other Details1", move the selection/active cell to D7
# - Clear the value in "Share Other Details1"!B5 (was 86)

$wb = $excel.ActiveWorkbook

$wsShare1 = $wb.Worksheets.Item("Share Other Details1")

# Clear out the stray value left in B5
$wsShare1.Range("B5").Value = $null

# Activate "Share Other Details1" and set its selection to D7,
# which makes it the active/selected tab (activeTab moves 2 -> 3)
# and leaves "Modify Transaction" no longer tab-selected.
$wsShare1.Activate()
$wsShare1.Range("D7").Select()
